# Final Design Phase Update
# - Change the highlight colour of three list-item headings from yellow to
#   magenta (both the run text and the paragraph-mark run properties).
# - Move the "_GoBack" bookmark from the "Website Screenshots" paragraph to
#   the "Prototype Screenshots" paragraph.

$d = $word.ActiveDocument

# wdPink (5) is serialised as w:highlight w:val="magenta" in OOXML.
$wdPink = 5

$headings = @("Prototype Screenshots", "File & Record Definitions", "Test Data")

foreach ($heading in $headings) {
    $rng = $d.Content
    $found = $rng.Find.Execute($heading, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        $para.Range.Font.HighlightColorIndex = $wdPink
    }
}

# Move the _GoBack bookmark onto the "Prototype Screenshots" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng2 = $d.Content
$rng2.Find.Execute("Prototype Screenshots", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($rng2.Start, $rng2.Start)
$d.Bookmarks.Add("_GoBack", $target)
